$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: WaitingTime column is inserted before NumberMeasuringFields,
# shifting L / A / B one column to the right (H..L rotate right by one).
$ws.Range("H1").Value = "WaitingTime"
$ws.Range("I1").Value = "NumberMeasuringFields"
$ws.Range("J1").Value = "L"
$ws.Range("K1").Value = "A"
$ws.Range("L1").Value = "B"

# Row 3: French model/description row for columns A-G (enum/description row)
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"

# Columns H-L of row 3 hold an empty (but present) text cell, matching
# the template's trailing columns. A bare "" clears the cell entirely, so
# force a text cell via the quote-prefix trick and then drop the resulting
# formatting so the cell keeps the default style.
$ws.Range("H3").Value = "'"
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "'"
$ws.Range("K3").Value = "'"
$ws.Range("L3").Value = "'"
$ws.Range("H3:L3").ClearFormats()
